# Update "想去人数" (want-to-go count) values for a few events that are
# duplicated across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12
$ws1.Range("F4").Value = 3329
$ws1.Range("F5").Value = 159
$ws1.Range("F7").Value = 160

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 12
$ws4.Range("F8").Value = 3329
$ws4.Range("F9").Value = 159
$ws4.Range("F12").Value = 160
